# Fix order of certificates on Sheet2 (A4:B8).
# The last two certificate rows (SEO, PR) are moved up to the top of the
# block, and the remaining rows (MySQL, AdvSQL, PowerBi) shift down to
# follow them - i.e. rows 4-8 are rotated left by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Capture current values for the block A4:B8 before overwriting anything.
$rows = @()
for ($r = 4; $r -le 8; $r++) {
    $rows += ,@($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 2).Value2)
}

# New order: rows [7,8,4,5,6] (1-based source rows) -> rows [4,5,6,7,8]
$newOrder = @(3, 4, 0, 1, 2)  # zero-based indices into $rows, matching old rows 7,8,4,5,6

$destRow = 4
foreach ($idx in $newOrder) {
    $pair = $rows[$idx]
    $ws.Cells.Item($destRow, 1).Value = $pair[0]
    $ws.Cells.Item($destRow, 2).Value = $pair[1]
    $destRow++
}
